$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-12-08 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-09 Monday", 2) | Out-Null

# Update each equation cell in the table, in row-major order, to match the
# source document exactly (avoids ambiguity from duplicate equation values).
$t = $d.Tables.Item(1)
$values = @(
    "2+40=42",
    "91-11=80",
    "19+54=73",
    "7-0=7",
    "69-57=12",
    "20+3=23",
    "12-4=8",
    "40-28=12",
    "26+1=27",
    "55-32=23",
    "54-1=53",
    "30+44=74",
    "70+16=86",
    "17+57=74",
    "62+28=90",
    "94-15=79",
    "72-67=5",
    "28+44=72",
    "59-44=15",
    "86-61=25",
    "52+22=74",
    "51-38=13",
    "86-12=74",
    "7+3=10",
    "42+38=80",
    "92-72=20",
    "42-1=41",
    "23+35=58",
    "27+42=69",
    "14+80=94",
    "78-31=47",
    "93-72=21",
    "29+14=43",
    "36-18=18",
    "53-12=41",
    "2+23=25",
    "61-25=36",
    "41+0=41",
    "54-33=21",
    "52-8=44",
    "85-54=31",
    "39+3=42",
    "13+23=36",
    "22+0=22",
    "64-0=64",
    "8+11=19",
    "36+50=86",
    "40+49=89",
    "58+9=67",
    "60+18=78",
    "16-11=5",
    "56-37=19",
    "36-1=35",
    "56-33=23",
    "34-7=27",
    "10+13=23",
    "74-14=60",
    "70-29=41",
    "36+14=50",
    "10+56=66",
    "78-59=19",
    "56-11=45",
    "17-16=1",
    "80-64=16",
    "73-6=67",
    "94-52=42",
    "10+58=68",
    "29+60=89",
    "43+54=97",
    "2+40=42",
    "43-1=42",
    "39-18=21",
    "63-36=27",
    "3-0=3",
    "85-65=20",
    "46+10=56",
    "30+57=87",
    "80-6=74",
    "28+57=85",
    "1+35=36",
    "47-5=42",
    "23+57=80",
    "84+11=95",
    "27+48=75",
    "8+0=8",
    "64-9=55",
    "10+54=64",
    "28-12=16",
    "38+2=40",
    "44+34=78",
    "97-3=94",
    "32+46=78",
    "33+27=60",
    "84-39=45",
    "39+31=70",
    "13+27=40",
    "9+80=89",
    "30+42=72",
    "33-21=12",
    "3+37=40"
)

$rows = 20
$cols = 5
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Output "Done: updated $idx cells"
